$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Update F8: "Por iniciar" -> "En proceso"
$ws.Range("F8").Value = "En proceso"

# Update N8: blank -> 2 (this cascades through the shared formulas in row 8)
$ws.Range("N8").Value = 2

# Re-merge these cells so they are re-appended at the end of the mergeCells list
$ws.Range("AZ4:BA4").UnMerge()
$ws.Range("AO4:AP4").UnMerge()
$ws.Range("AR4:AS4").UnMerge()
$ws.Range("AU4:AV4").UnMerge()
$ws.Range("AX4:AY4").UnMerge()
$ws.Range("AZ4:BA4").Merge()
$ws.Range("AO4:AP4").Merge()
$ws.Range("AR4:AS4").Merge()
$ws.Range("AU4:AV4").Merge()
$ws.Range("AX4:AY4").Merge()

# Update active cell selection to N8
$ws.Range("N8").Select()
